$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$r2A = @'
This test is to get description components.

'@
$ws.Range("A2").Value = $r2A
$ws.Range("B2").NumberFormat = "@"
$r2B = @'
0
'@
$ws.Range("B2").Value = $r2B
$r2C = @'
<class 'AutomationFramework.page_objects.hardware.hardware.Hardware'>
'@
$ws.Range("C2").Value = $r2C
$r2D = @'
hw_component_description
'@
$ws.Range("D2").Value = $r2D
$r2G = @'
<get>
  <filter>
    <components xmlns="http://openconfig.net/yang/platform">
      <component>
        <name>CHASIS</name>
        <state>
          <description></description>
        </state>
      </component>
    </components>
  </filter>
</get>
'@
$ws.Range("G2").Value = $r2G
$r2J = @'
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:a9db80e4-4328-4d02-8edf-2b21b768cab6"
 xmlns:ncx="http://netconfcentral.org/ns/yuma-ncx"
 ncx:last-modified="2020-10-07T13:51:28Z" ncx:etag="814"
 xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
 <data></data>
</rpc-reply>
'@
$ws.Range("J2").Value = $r2J

# Row 3
$r3A = @'
This test is to get the hardware version components.

'@
$ws.Range("A3").Value = $r3A
$ws.Range("B3").NumberFormat = "@"
$r3B = @'
0
'@
$ws.Range("B3").Value = $r3B
$r3C = @'
<class 'AutomationFramework.page_objects.hardware.hardware.Hardware'>
'@
$ws.Range("C3").Value = $r3C
$r3D = @'
hw_component_hardware_version
'@
$ws.Range("D3").Value = $r3D
$r3G = @'
<get>
  <filter>
    <components xmlns="http://openconfig.net/yang/platform">
      <component>
        <name>CHASIS</name>
        <state>
          <hardware-version></hardware-version>
        </state>
      </component>
    </components>
  </filter>
</get>
'@
$ws.Range("G3").Value = $r3G
$r3J = @'
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:770b8d52-6123-4d3f-9c3b-94a954bd3ed3"
 xmlns:ncx="http://netconfcentral.org/ns/yuma-ncx"
 ncx:last-modified="2020-10-07T13:51:28Z" ncx:etag="814"
 xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
 <data></data>
</rpc-reply>
'@
$ws.Range("J3").Value = $r3J

# Row 4
$r4A = @'
This test is to get the component id.

'@
$ws.Range("A4").Value = $r4A
$ws.Range("B4").NumberFormat = "@"
$r4B = @'
0
'@
$ws.Range("B4").Value = $r4B
$r4C = @'
<class 'AutomationFramework.page_objects.hardware.hardware.Hardware'>
'@
$ws.Range("C4").Value = $r4C
$r4D = @'
hw_component_id
'@
$ws.Range("D4").Value = $r4D
$r4G = @'
<get>
  <filter>
    <components xmlns="http://openconfig.net/yang/platform">
      <component>
        <name>CHASIS</name>
        <state>
          <id></id>
        </state>
      </component>
    </components>
  </filter>
</get>
'@
$ws.Range("G4").Value = $r4G
$r4J = @'
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:e182e267-14e7-4605-924e-8b7fc69051a4"
 xmlns:ncx="http://netconfcentral.org/ns/yuma-ncx"
 ncx:last-modified="2020-10-07T13:51:28Z" ncx:etag="814"
 xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
 <data></data>
</rpc-reply>
'@
$ws.Range("J4").Value = $r4J

# Row 5
$r5A = @'
This test is to get location components.

'@
$ws.Range("A5").Value = $r5A
$ws.Range("B5").NumberFormat = "@"
$r5B = @'
0
'@
$ws.Range("B5").Value = $r5B
$r5C = @'
<class 'AutomationFramework.page_objects.hardware.hardware.Hardware'>
'@
$ws.Range("C5").Value = $r5C
$r5D = @'
hw_component_location
'@
$ws.Range("D5").Value = $r5D
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$r5G = @'
<get>
  <filter>
    <components xmlns="http://openconfig.net/yang/platform">
      <component>
        <name>CHASIS</name>
        <state>
          <location></location>
        </state>
      </component>
    </components>
  </filter>
</get>
'@
$ws.Range("G5").Value = $r5G
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""
$r5J = @'
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:ced8c9f7-efd4-46a8-b0a2-e4fab569c00c"
 xmlns:ncx="http://netconfcentral.org/ns/yuma-ncx"
 ncx:last-modified="2020-10-07T13:51:28Z" ncx:etag="814"
 xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
 <data></data>
</rpc-reply>
'@
$ws.Range("J5").Value = $r5J
